# Updated cryptos list (price + 1h volume change) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $range = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the value as literal text,
    # even for strings that look numeric (e.g. "0.999", "68.690.28").
    $range.Formula = "`'" + $text
    # Reset to the default style so no number-format/quote-prefix style
    # gets attached to the cell (matches the original, style-less cells).
    $range.Style = "Normal"
}

Set-TextCell "D2" "68.690.28"
Set-TextCell "E2" "  +1.35%  "
Set-TextCell "D3" "3.730.14"
Set-TextCell "E3" "  -2.11%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.11%  "
Set-TextCell "D5" "601.67"
Set-TextCell "E5" "  -0.23%  "
Set-TextCell "D6" "169.36"
Set-TextCell "E6" "  -2.08%  "
Set-TextCell "D7" "3.720.69"
Set-TextCell "E7" "  -2.30%  "
Set-TextCell "E8" "  -0.03%  "
Set-TextCell "D9" "0.535"
Set-TextCell "E9" "  +0.59%  "
Set-TextCell "D10" "0.164"
Set-TextCell "D11" "6.35"
Set-TextCell "E11" "  +1.09%  "
Set-TextCell "D12" "0.461"
Set-TextCell "E12" "  -1.15%  "
Set-TextCell "D13" "38.15"
Set-TextCell "E13" "  -1.54%  "
Set-TextCell "D14" "0.0000245"
Set-TextCell "E14" "  +0.34%  "
Set-TextCell "D15" "4.342.35"
Set-TextCell "E15" "  -2.23%  "
Set-TextCell "D16" "3.715.16"
Set-TextCell "E16" "  -2.49%  "
Set-TextCell "D17" "68.619.41"
Set-TextCell "E17" "  +1.23%  "
Set-TextCell "D18" "7.29"
Set-TextCell "E18" "  +0.38%  "
Set-TextCell "E19" "  +0.99%  "
Set-TextCell "D20" "17.15"
Set-TextCell "E20" "  +1.18%  "
Set-TextCell "D21" "496.31"
Set-TextCell "E21" "  +0.58%  "
Set-TextCell "D22" "10.37"
Set-TextCell "E22" "  +12.08%  "
Set-TextCell "D23" "0.725"
Set-TextCell "E23" "  -2.93%  "
Set-TextCell "D24" "85.19"
Set-TextCell "E24" "  -1.09%  "
Set-TextCell "D25" "0.0000144"
Set-TextCell "E25" "  -2.11%  "
Set-TextCell "D26" "2.31"
Set-TextCell "E26" "  -2.84%  "
Set-TextCell "D27" "12.46"
Set-TextCell "E27" "  +1.03%  "
Set-TextCell "D28" "10.16"
Set-TextCell "E28" "  -0.98%  "
Set-TextCell "E29" "  -0.09%  "
Set-TextCell "D30" "2.58"
Set-TextCell "E30" "  +5.48%  "
Set-TextCell "D31" "2.97"
Set-TextCell "E31" "  -0.82%  "
Set-TextCell "D32" "7.98"
Set-TextCell "E32" "  +2.45%  "
Set-TextCell "D33" "31.63"
Set-TextCell "E33" "  -4.79%  "
Set-TextCell "D34" "3.865.35"
Set-TextCell "E34" "  -2.24%  "
Set-TextCell "E35" "  -1.10%  "
Set-TextCell "D36" "3.655.19"
Set-TextCell "E36" "  -2.42%  "
Set-TextCell "D37" "0.998"
Set-TextCell "E37" "  -0.17%  "
Set-TextCell "E38" "  -0.48%  "
Set-TextCell "D39" "5.85"
Set-TextCell "E39" "  +0.23%  "
Set-TextCell "E40" "  -0.88%  "
Set-TextCell "D41" "0.326"
Set-TextCell "E41" "  -1.21%  "
Set-TextCell "D42" "438.66"
Set-TextCell "E42" "  -5.34%  "
Set-TextCell "D43" "48.81"
Set-TextCell "E43" "  -0.62%  "
Set-TextCell "D44" "1.98"
Set-TextCell "E44" "  -1.69%  "
Set-TextCell "D45" "2.89"
Set-TextCell "E45" "  +0.76%  "
Set-TextCell "E46" "  +0.52%  "
Set-TextCell "E47" "  +0.00%  "
Set-TextCell "D48" "40.67"
Set-TextCell "E48" "  -1.38%  "
Set-TextCell "D49" "141.13"
Set-TextCell "E49" "  +1.52%  "
Set-TextCell "D50" "0.0354"
Set-TextCell "E50" "  +0.77%  "
Set-TextCell "D51" "2.761.04"
Set-TextCell "E51" "  -3.22%  "
